$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Task-row value edits (rows 9-33 "Daily Time Sheet" detail grid) ---
$ws.Range("L14").Value = 1
$ws.Range("O14").Value = 1

$ws.Range("O19").Value = 1.5

$ws.Range("Q24").Value = 1
$ws.Range("Q25").Value = 1
$ws.Range("Q26").Value = 1
$ws.Range("Q27").Value = 2
$ws.Range("Q28").Value = 2
$ws.Range("Q29").Value = 1
$ws.Range("Q30").Value = 1
$ws.Range("Q31").Value = 1

$ws.Range("R32").Value = 3
$ws.Range("R33").Value = 1

# --- View state: scroll back to top-left, change zoom, move selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$win.Zoom = 100
$ws.Range("S15").Select()
